# Last update from John!
# Applies the ISR Serial ICD sheet updates: renamed/added lookup values and
# highlights the newly-added mag/accel calibration rows in bold blue.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Text content updates -------------------------------------------------

# "minAcro" -> "Misc Values" (row 5, columns B and D both reference the same
# lookup value)
$ws.Range("B5").Value = "Misc Values"
$ws.Range("D5").Value = "Misc Values"

# Row 18 gains new lookup entries for the mag calibration commands
$ws.Range("D18").Value = "Send mag cal values"
$ws.Range("B18").Value = "Write mag calibration values"

# New entries appended to the D column lookup list
$ws.Range("D20").Value = "Send Pressure Altitude"
$ws.Range("D24").Value = "Send Command in Detent Discretes"

# Row 16: D column changes from "Send calibration values" to the more
# specific "Send accel calibration values"
$ws.Range("D16").Value = "Send accel calibration values"

$ws.Range("D26").Value = "Send 100 Hz loop time"

# --- Formatting updates -----------------------------------------------------
# Highlight the accel/mag calibration rows (16-18, columns B & D) in bold
# blue to call out the newly added commands.

foreach ($addr in @("B16", "D16", "D17", "B18", "D18")) {
    $cell = $ws.Range($addr)
    $cell.Font.Bold = $true
    $cell.Font.Color = 15773696
}

# --- Selection -------------------------------------------------------------
$ws.Range("I30").Select()
